$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# Update existing cell C112: 6 -> 7
$ws.Range("C112").Value = 7

# Prepare formatting for the two new rows by copying the format of row 112
$ws.Range("A112:C112").Copy()
$ws.Range("A113:C114").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 113
$ws.Range("A113").Value = 44854
$ws.Range("B113").Value = "Tablas de edad, educion y sexo por bases de datos"
$ws.Range("C113").Value = 2

# New row 114
$ws.Range("A114").Value = 44855
$ws.Range("B114").Value = "Ear dream y reunion con emilse"
$ws.Range("C114").Value = 4

# Update selection to match the new active cell
$ws.Range("C114").Select()
